# Insert a new data row at sheet row 675 (2026/01/20, 09:00 reading) and
# push the existing 2026/12/29 .. 2027/01/05 rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 675..716 down to 676..717, opening up a blank row at 675.
$ws.Rows("675:675").Insert()

# A675 text looks like a date ("2026/01/20"); Excel would normally parse it
# into a date serial. Lead with an apostrophe to force plain text entry
# (same as typing '2026/01/20 into the cell), then strip the leftover
# "quote prefix" formatting so the cell keeps the sheet's default style.
$ws.Range("A675").Value = "'2026/01/20"
$ws.Range("A675").ClearFormats()

$ws.Range("B675").Value = "火"
$ws.Range("C675").Value = 9
$ws.Range("D675").Value = 201
